$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written with a leading apostrophe to force text entry
# (matches the source inline-string cells, avoids numeric auto-coercion
# for values like "483.65"), then the style is reset to Normal so the
# quote-prefix flag does not leave a residual style/format change.

$ws.Range("D2").Value = "'53.584.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.223.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.60%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'483.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.90%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'125.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.92%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  -5.42%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.230.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.63%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -7.36%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.02%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.67%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.619.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.51%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'21.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'53.487.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.229.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'Polkadot"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'3.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.21%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Chainlink"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'9.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'297.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'63.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.36%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.81%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.143"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'169.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.27%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0682"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.98%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -3.94%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.62%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +4.74%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -6.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'35.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.36%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.363"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.09%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.90%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'122.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -6.59%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'4.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.32%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.57%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.77%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'231.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.52%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.81%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'16.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.43%  "
$ws.Range("E51").Style = "Normal"
